# Updated cryptos list (Price / Volume(1h) columns) with latest scraped values.
# D-column "price" cells that look numeric (single decimal point, e.g. "579.61")
# are forced to text via NumberFormat "@" before assignment so Excel does not
# silently convert them to floating point numbers (which would both change the
# cell type away from text and introduce binary floating-point rounding, e.g.
# "579.61" -> 579.61000000000001). The Style is reset to "Normal" afterwards
# so no stray custom cell style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.644.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.446.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.49%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.034.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("E14").Value = "  -5.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.438.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.692.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.581.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.482.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.573.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("E51").Value = "  -0.08%  "

Write-Host "Updated cryptos list with latest price/volume data"
